# Apply updated profit figures across all sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 9902.75
$ws.Range("I46").Value = 6000
$ws.Range("J46").Value = 11203.667
$ws.Range("K46").Value = 18000
$ws.Range("L46").Value = 33611.001
$ws.Range("M46").Value = -17881
$ws.Range("N46").Value = -33849.001
$ws.Range("H60").Value = 9902.75
$ws.Range("I60").Value = 6000
$ws.Range("J60").Value = 11203.667
$ws.Range("K60").Value = 18000
$ws.Range("L60").Value = 33611.001
$ws.Range("M60").Value = -17516
$ws.Range("N60").Value = -34579.001
$ws.Range("H62").Value = 4067
$ws.Range("J62").Value = 4067
$ws.Range("L62").Value = 4067
$ws.Range("N62").Value = -5315
$ws.Range("H65").Value = 4067
$ws.Range("J65").Value = 4067
$ws.Range("L65").Value = 20335
$ws.Range("N65").Value = -26575
$ws.Range("H69").Value = 6333.3335
$ws.Range("I69").Value = 7000
$ws.Range("J69").Value = 5000
$ws.Range("K69").Value = 21000
$ws.Range("L69").Value = 15000
$ws.Range("M69").Value = -20126
$ws.Range("N69").Value = -16748
$ws.Range("H72").Value = 6333.3335
$ws.Range("I72").Value = 7000
$ws.Range("J72").Value = 5000
$ws.Range("K72").Value = 63000
$ws.Range("L72").Value = 45000
$ws.Range("M72").Value = -58632
$ws.Range("N72").Value = -53736
$ws.Range("H98").Value = 2092.2222
$ws.Range("I98").Value = 2092.2222
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2092.2222
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -594.2222000000002
$ws.Range("H112").Value = 4147.5
$ws.Range("J112").Value = 4175.185
$ws.Range("L112").Value = 12525.555
$ws.Range("N112").Value = -14741.555
$ws.Range("H113").Value = 3396.2593
$ws.Range("J113").Value = 3624.9285
$ws.Range("L113").Value = 3624.9285
$ws.Range("N113").Value = -10132.9285
$ws.Range("H116").Value = 9860209
$ws.Range("I116").Value = 27102080
$ws.Range("J116").Value = 7710.5713
$ws.Range("K116").Value = 27102080
$ws.Range("L116").Value = 7710.5713
$ws.Range("M116").Value = -27098638
$ws.Range("N116").Value = -14594.5713
$ws.Range("H122").Value = 2092.2222
$ws.Range("I122").Value = 2092.2222
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6276.6666
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3826.6666
$ws.Range("H131").Value = 1230.2727
$ws.Range("I131").Value = 1273.8
$ws.Range("K131").Value = 3821.4
$ws.Range("M131").Value = 1218.6
$ws.Range("H132").Value = 12565.642
$ws.Range("I132").Value = 2668.6
$ws.Range("J132").Value = 14867.279
$ws.Range("K132").Value = 8005.799999999999
$ws.Range("L132").Value = 44601.837
$ws.Range("M132").Value = -5475.799999999999
$ws.Range("N132").Value = -49661.837
$ws.Range("H138").Value = 6755.7617
$ws.Range("I138").Value = 1539.875
$ws.Range("J138").Value = 7983.0293
$ws.Range("K138").Value = 4619.625
$ws.Range("L138").Value = 23949.0879
$ws.Range("M138").Value = 520.375
$ws.Range("N138").Value = -34229.0879
$ws.Range("N98").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1746921.8
$ws.Range("I2").Value = 2910042
$ws.Range("K2").Value = 2910042
$ws.Range("M2").Value = -2909929
$ws.Range("H36").Value = 6242
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("H102").Value = 653843.3
$ws.Range("I102").Value = 1958466.9
$ws.Range("K102").Value = 1958466.9
$ws.Range("M102").Value = -1956844.9
$ws.Range("H116").Value = 1746921.8
$ws.Range("I116").Value = 2910042
$ws.Range("K116").Value = 2910042
$ws.Range("M116").Value = -2907748
$ws.Range("N36").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1746921.8
$ws.Range("I3").Value = 2910042
$ws.Range("K3").Value = 2910042
$ws.Range("M3").Value = -2909928
$ws.Range("H134").Value = 2074.4219
$ws.Range("J134").Value = 3949.3
$ws.Range("L134").Value = 11847.9
$ws.Range("N134").Value = -16917.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 40000
$ws.Range("J28").Value = 40000
$ws.Range("L28").Value = 40000
$ws.Range("N28").Value = -40490
$ws.Range("H31").Value = 2308
$ws.Range("I31").Value = 1170.2778
$ws.Range("J31").Value = 3770.7856
$ws.Range("K31").Value = 1170.2778
$ws.Range("L31").Value = 3770.7856
$ws.Range("M31").Value = -875.2778000000001
$ws.Range("N31").Value = -4360.7856
$ws.Range("H34").Value = 2308
$ws.Range("I34").Value = 1170.2778
$ws.Range("J34").Value = 3770.7856
$ws.Range("K34").Value = 1170.2778
$ws.Range("L34").Value = 3770.7856
$ws.Range("M34").Value = -968.2778000000001
$ws.Range("N34").Value = -4174.7856
$ws.Range("H62").Value = 30470.908
$ws.Range("I62").Value = 4168.3335
$ws.Range("J62").Value = 40334.375
$ws.Range("K62").Value = 4168.3335
$ws.Range("L62").Value = 40334.375
$ws.Range("M62").Value = -3544.3335
$ws.Range("N62").Value = -41582.375
$ws.Range("H65").Value = 30470.908
$ws.Range("I65").Value = 4168.3335
$ws.Range("J65").Value = 40334.375
$ws.Range("K65").Value = 20841.6675
$ws.Range("L65").Value = 201671.875
$ws.Range("M65").Value = -17721.6675
$ws.Range("N65").Value = -207911.875
$ws.Range("H134").Value = 2641.6765
$ws.Range("I134").Value = 2588.0625
$ws.Range("J134").Value = 3499.5
$ws.Range("K134").Value = 7764.1875
$ws.Range("L134").Value = 10498.5
$ws.Range("M134").Value = -5229.1875
$ws.Range("N134").Value = -15568.5
$ws.Range("H141").Value = 113721.336
$ws.Range("I141").Value = 39498
$ws.Range("K141").Value = 39498
$ws.Range("M141").Value = -34318

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 17377.6
$ws.Range("J125").Value = 19999
$ws.Range("L125").Value = 59997
$ws.Range("N125").Value = -69837
$ws.Range("H129").Value = 2998
$ws.Range("I129").Value = 2164
$ws.Range("K129").Value = 6492
$ws.Range("M129").Value = -1492
$ws.Range("H137").Value = 5800653.5
$ws.Range("I137").Value = 4975
$ws.Range("J137").Value = 6887343.5
$ws.Range("K137").Value = 14925
$ws.Range("L137").Value = 20662030.5
$ws.Range("M137").Value = -9825
$ws.Range("N137").Value = -20672230.5
$ws.Range("H138").Value = 3158.9443
$ws.Range("I138").Value = 2991.1333
$ws.Range("J138").Value = 3998
$ws.Range("K138").Value = 8973.3999
$ws.Range("L138").Value = 11994
$ws.Range("M138").Value = -3833.3999
$ws.Range("N138").Value = -22274
$ws.Range("H139").Value = 2056.5
$ws.Range("I139").Value = 2056.5
$ws.Range("K139").Value = 6169.5
$ws.Range("M139").Value = -1029.5
$ws.Range("H141").Value = 8693.789000000001
$ws.Range("I141").Value = 5938.3076
$ws.Range("J141").Value = 14664
$ws.Range("K141").Value = 17814.9228
$ws.Range("L141").Value = 43992
$ws.Range("M141").Value = -12634.9228
$ws.Range("N141").Value = -54352

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 9286.1
$ws.Range("J99").Value = 27499.5
$ws.Range("L99").Value = 27499.5
$ws.Range("N99").Value = -31991.5
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("H113").Value = 3408.25
$ws.Range("I113").Value = 3211
$ws.Range("K113").Value = 3211
$ws.Range("M113").Value = -1041
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 87
$ws.Range("I55").Value = 96.666664
$ws.Range("J55").Value = 77.333336
$ws.Range("K55").Value = 96.666664
$ws.Range("L55").Value = 77.333336
$ws.Range("M55").Value = 76.333336
$ws.Range("N55").Value = -423.333336
$ws.Range("H99").Value = 57999.8
$ws.Range("I99").Value = 24999.5
$ws.Range("K99").Value = 24999.5
$ws.Range("M99").Value = -22004.5
$ws.Range("H122").Value = 20412734
$ws.Range("I122").Value = 2666.3333
$ws.Range("J122").Value = 35720284
$ws.Range("K122").Value = 7998.999899999999
$ws.Range("L122").Value = 107160852
$ws.Range("M122").Value = -5548.999899999999
$ws.Range("N122").Value = -107165752

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2432.0715
$ws.Range("I107").Value = 2691.9565
$ws.Range("K107").Value = 8075.869499999999
$ws.Range("M107").Value = -6155.869499999999
$ws.Range("H108").Value = 87969
$ws.Range("J108").Value = 87969
$ws.Range("L108").Value = 87969
$ws.Range("N108").Value = -95649
